{"js": "// Replace the date line and every three-digit-division answer cell with\n// the newly generated values. Each old string is unique in the document,\n// so a simple search-and-replace per pair is unambiguous.\nconst replacements = [\n  [\"2025-02-01 Saturday\", \"2025-02-02 Sunday\"],\n  [\"575\u00f79=63, 8\", \"490\u00f74=122, 2\"],\n  [\"619\u00f72=309, 1\", \"191\u00f76=31, 5\"],\n  [\"746\u00f72=373, 0\", \"769\u00f78=96, 1\"],\n  [\"241\u00f75=48, 1\", \"183\u00f76=30, 3\"],\n  [\"615\u00f77=87, 6\", \"944\u00f78=118, 0\"],\n  [\"616\u00f73=205, 1\", \"295\u00f72=147, 1\"],\n  [\"398\u00f78=49, 6\", \"905\u00f74=226, 1\"],\n  [\"910\u00f73=303, 1\", \"496\u00f74=124, 0\"],\n  [\"216\u00f79=24, 0\", \"242\u00f74=60, 2\"],\n  [\"168\u00f76=28, 0\", \"741\u00f76=123, 3\"],\n  [\"196\u00f76=32, 4\", \"276\u00f74=69, 0\"],\n  [\"409\u00f78=51, 1\", \"525\u00f78=65, 5\"],\n  [\"569\u00f73=189, 2\", \"254\u00f79=28, 2\"],\n  [\"461\u00f73=153, 2\", \"997\u00f75=199, 2\"],\n  [\"659\u00f79=73, 2\", \"147\u00f78=18, 3\"],\n  [\"235\u00f72=117, 1\", \"460\u00f78=57, 4\"],\n  [\"313\u00f76=52, 1\", \"821\u00f75=164, 1\"],\n  [\"408\u00f72=204, 0\", \"839\u00f73=279, 2\"],\n  [\"847\u00f75=169, 2\", \"216\u00f73=72, 0\"],\n  [\"591\u00f75=118, 1\", \"647\u00f73=215, 2\"],\n  [\"710\u00f77=101, 3\", \"707\u00f77=101, 0\"],\n  [\"927\u00f77=132, 3\", \"310\u00f73=103, 1\"],\n  [\"865\u00f79=96, 1\", \"341\u00f78=42, 5\"],\n  [\"184\u00f79=20, 4\", \"113\u00f78=14, 1\"],\n  [\"122\u00f78=15, 2\", \"269\u00f72=134, 1\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and every three-digit-division answer cell with\n# the newly generated values. Each old string is unique in the document,\n# so Find/Replace (wdReplaceAll) per pair is unambiguous and safe.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-02-01 Saturday\", \"2025-02-02 Sunday\"),\n    @(\"575\u00f79=63, 8\", \"490\u00f74=122, 2\"),\n    @(\"619\u00f72=309, 1\", \"191\u00f76=31, 5\"),\n    @(\"746\u00f72=373, 0\", \"769\u00f78=96, 1\"),\n    @(\"241\u00f75=48, 1\", \"183\u00f76=30, 3\"),\n    @(\"615\u00f77=87, 6\", \"944\u00f78=118, 0\"),\n    @(\"616\u00f73=205, 1\", \"295\u00f72=147, 1\"),\n    @(\"398\u00f78=49, 6\", \"905\u00f74=226, 1\"),\n    @(\"910\u00f73=303, 1\", \"496\u00f74=124, 0\"),\n    @(\"216\u00f79=24, 0\", \"242\u00f74=60, 2\"),\n    @(\"168\u00f76=28, 0\", \"741\u00f76=123, 3\"),\n    @(\"196\u00f76=32, 4\", \"276\u00f74=69, 0\"),\n    @(\"409\u00f78=51, 1\", \"525\u00f78=65, 5\"),\n    @(\"569\u00f73=189, 2\", \"254\u00f79=28, 2\"),\n    @(\"461\u00f73=153, 2\", \"997\u00f75=199, 2\"),\n    @(\"659\u00f79=73, 2\", \"147\u00f78=18, 3\"),\n    @(\"235\u00f72=117, 1\", \"460\u00f78=57, 4\"),\n    @(\"313\u00f76=52, 1\", \"821\u00f75=164, 1\"),\n    @(\"408\u00f72=204, 0\", \"839\u00f73=279, 2\"),\n    @(\"847\u00f75=169, 2\", \"216\u00f73=72, 0\"),\n    @(\"591\u00f75=118, 1\", \"647\u00f73=215, 2\"),\n    @(\"710\u00f77=101, 3\", \"707\u00f77=101, 0\"),\n    @(\"927\u00f77=132, 3\", \"310\u00f73=103, 1\"),\n    @(\"865\u00f79=96, 1\", \"341\u00f78=42, 5\"),\n    @(\"184\u00f79=20, 4\", \"113\u00f78=14, 1\"),\n    @(\"122\u00f78=15, 2\", \"269\u00f72=134, 1\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
